# Update "想去人数" (interested-count) figures across the four sheets.
# Mapping derived from the OOXML diff: sheet name -> { cell -> new value }

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F2"  = 1553
        "F5"  = 8558
        "F7"  = 94
        "F8"  = 1231
        "F9"  = 27
        "F10" = 233
        "F13" = 93
        "F14" = 271
        "F17" = 1376
        "F18" = 1301
        "F19" = 561
        "F21" = 1316
        "F23" = 200
        "F26" = 53
        "F27" = 263
        "F28" = 1055
        "F31" = 194
        "F32" = 161
        "F35" = 591
        "F37" = 110
        "F39" = 135
        "F43" = 193
    }
    "演出" = @{
        "F2"  = 39
        "F5"  = 22
        "F8"  = 7
        "F9"  = 22
        "F12" = 218
        "F17" = 656
        "F21" = 53
        "F24" = 915
        "F26" = 1020
        "F27" = 100
        "F28" = 631
        "F35" = 154
        "F37" = 19
    }
    "本地生活" = @{
        "F5"  = 871
        "F9"  = 1918
        "F10" = 2906
    }
    "全部类型" = @{
        "F4"  = 22
        "F6"  = 722
        "F8"  = 8558
        "F11" = 7
        "F12" = 1918
        "F13" = 2906
        "F16" = 218
        "F17" = 94
        "F18" = 1231
        "F22" = 93
        "F23" = 271
        "F24" = 1376
        "F25" = 1301
        "F26" = 1316
        "F27" = 200
        "F28" = 53
        "F29" = 263
        "F32" = 53
        "F34" = 915
        "F35" = 194
        "F37" = 161
        "F39" = 591
        "F40" = 631
        "F43" = 154
        "F44" = 193
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
